$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.371.89"
$ws.Range("E2").Value = "  +0.44%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.227.96"
$ws.Range("E3").Value = "  -0.55%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.14%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'244.86"
$ws.Range("E5").Value = "  -0.46%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.52%  "

# Row 7
$ws.Range("D7").Value = "'74.32"
$ws.Range("E7").Value = "  -2.40%  "

# Row 8
$ws.Range("E8").Value = "  +0.20%  "

# Row 9
$ws.Range("E9").Value = "  -0.53%  "

# Row 10
$ws.Range("D10").Value = "'43.21"
$ws.Range("E10").Value = "  +4.08%  "

# Row 11
$ws.Range("E11").Value = "  +1.77%  "

# Row 12
$ws.Range("D12").Value = "'7.13"
$ws.Range("E12").Value = "  +1.02%  "

# Row 13
$ws.Range("D13").Value = "'0.103"
$ws.Range("E13").Value = "  +0.20%  "

# Row 14
$ws.Range("D14").Value = "'14.44"
$ws.Range("E14").Value = "  -1.23%  "

# Row 15
$ws.Range("D15").Value = "'0.851"
$ws.Range("E15").Value = "  -0.22%  "

# Row 16
$ws.Range("D16").Value = "2.226.01"
$ws.Range("E16").Value = "  -0.40%  "

# Row 17
$ws.Range("D17").Value = "42.154.12"
$ws.Range("E17").Value = "  +0.45%  "

# Row 18
$ws.Range("E18").Value = "  +12.27%  "

# Row 19
$ws.Range("E19").Value = "  +1.76%  "

# Row 20
$ws.Range("D20").Value = "'72.08"
$ws.Range("E20").Value = "  +0.73%  "

# Row 21
$ws.Range("D21").Value = "'10.02"
$ws.Range("E21").Value = "  +35.53%  "

# Row 22
$ws.Range("D22").Value = "'230.98"
$ws.Range("E22").Value = "  +0.29%  "

# Row 23
$ws.Range("D23").Value = "'2.16"
$ws.Range("E23").Value = "  -5.23%  "

# Row 24
$ws.Range("D24").Value = "'11.81"
$ws.Range("E24").Value = "  +5.52%  "

# Row 25
$ws.Range("E25").Value = "  +0.06%  "

# Row 26
$ws.Range("E26").Value = "  -0.67%  "

# Row 27
$ws.Range("E27").Value = "  +1.10%  "

# Row 28
$ws.Range("E28").Value = "  +6.38%  "

# Row 29
$ws.Range("D29").Value = "'167.33"
$ws.Range("E29").Value = "  -0.86%  "

# Row 30
$ws.Range("D30").Value = "'21.08"
$ws.Range("E30").Value = "  +2.77%  "

# Row 31
$ws.Range("D31").Value = "'5.81"
$ws.Range("E31").Value = "  +17.88%  "

# Row 32
$ws.Range("E32").Value = "  -2.22%  "

# Row 33
$ws.Range("D33").Value = "'0.117"
$ws.Range("E33").Value = "  -1.19%  "

# Row 34 - was Stellar, now InjectiveProtocol (rows 34/35 swapped)
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "'29.48"
$ws.Range("E34").Value = "  -8.27%  "

# Row 35 - was InjectiveProtocol, now Stellar
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "'0.124"
$ws.Range("E35").Value = "  -0.65%  "

# Row 36
$ws.Range("D36").Value = "'4.41"
$ws.Range("E36").Value = "  -0.56%  "

# Row 37
$ws.Range("E37").Value = "  +2.38%  "

# Row 38
$ws.Range("E38").Value = "  -4.88%  "

# Row 39
$ws.Range("D39").Value = "'2.17"
$ws.Range("E39").Value = "  -0.01%  "

# Row 40
$ws.Range("D40").Value = "'5.62"
$ws.Range("E40").Value = "  -3.26%  "

# Row 41
$ws.Range("D41").Value = "'63.22"
$ws.Range("E41").Value = "  +4.57%  "

# Row 42
$ws.Range("E42").Value = "  -0.10%  "

# Row 43
$ws.Range("D43").Value = "'8.84"
$ws.Range("E43").Value = "  +2.19%  "

# Row 44
$ws.Range("D44").Value = "'104.74"
$ws.Range("E44").Value = "  -6.87%  "

# Row 45
$ws.Range("E45").Value = "  +3.12%  "

# Row 46
$ws.Range("D46").Value = "'0.994"
$ws.Range("E46").Value = "  -0.17%  "

# Row 47
$ws.Range("D47").Value = "'2.39"
$ws.Range("E47").Value = "  +6.24%  "

# Row 48
$ws.Range("E48").Value = "  +0.19%  "

# Row 49
$ws.Range("E49").Value = "  +1.91%  "

# Row 50
$ws.Range("E50").Value = "  +0.71%  "

# Row 51
$ws.Range("E51").Value = "  -0.84%  "
